$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (K = strikeouts), rows 2-21.
# Regenerated from Strike# to K, recalculated s_vals.
$gValues = @{
    2  = 0
    3  = 1
    4  = 3
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 3
    10 = 0
    11 = 2
    12 = 4
    13 = 1
    14 = 3
    15 = 5
    16 = 3
    17 = 7
    18 = 4
    19 = 1
    20 = 5
    21 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
